# Auto-Informes/plantillas/archivo.xlsx
# "Se incluye función de descargar informe"
#
# The sheet's header row used the misspelled "region" in G1; fix it to the
# correctly accented "región". This is the only real content change in the
# commit - it forces the shared-strings table to be rewritten (the old
# "region" entry at index 6 drops out and "región" is appended at the end),
# which is why every other shared-string index in the sheet shifts down by
# one even though no other cell's text actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "región"

# Leave the cursor on G2, matching the saved selection in the workbook.
$ws.Range("G2").Select()
